# Auto update Excel log: append newly-logged sensor readings to PIR, Humidity, and Temperature sheets
$wb = $excel.ActiveWorkbook

# --- PIR sheet: rows 247-257 ---
$ws = $wb.Worksheets.Item("PIR")

$newRows = @(
    ,@("2026-01-28", "16:30:49", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:51", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:53", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:55", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:59", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:31:05", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:31:10", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:31:15", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:31:20", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:31:25", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:31:30", "16:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 247
$endRow = 257
# Force text storage (matches source data: dates/percentages/etc. kept as literal strings)
$ws.Range("A$startRow`:F$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# --- Humidity sheet: rows 239-247 ---
$ws = $wb.Worksheets.Item("Humidity")

$newRows = @(
    ,@("2026-01-28", "16:30:48", "16:00", "Bathroom", "87.9%", "Active")
    ,@("2026-01-28", "16:30:52", "16:00", "Bathroom", "87.9%", "Active")
    ,@("2026-01-28", "16:30:54", "16:00", "Bathroom", "87.0%", "Active")
    ,@("2026-01-28", "16:31:00", "16:00", "Bathroom", "88.0%", "Active")
    ,@("2026-01-28", "16:31:04", "16:00", "Bathroom", "87.1%", "Active")
    ,@("2026-01-28", "16:31:16", "16:00", "Bathroom", "87.1%", "Active")
    ,@("2026-01-28", "16:31:20", "16:00", "Bathroom", "88.0%", "Active")
    ,@("2026-01-28", "16:31:24", "16:00", "Bathroom", "87.1%", "Active")
    ,@("2026-01-28", "16:31:28", "16:00", "Bathroom", "88.1%", "Active")
)

$startRow = 239
$endRow = 247
# Force text storage (matches source data: dates/percentages/etc. kept as literal strings)
$ws.Range("A$startRow`:F$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# --- Temperature sheet: rows 240-248 ---
$ws = $wb.Worksheets.Item("Temperature")

$newRows = @(
    ,@("2026-01-28", "16:30:48", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:30:53", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:30:55", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:31:01", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:31:05", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:31:17", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:31:21", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:31:25", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:31:29", "16:00", "Bathroom", "22.8C", "Active")
)

$startRow = 240
$endRow = 248
# Force text storage (matches source data: dates/percentages/etc. kept as literal strings)
$ws.Range("A$startRow`:F$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}
